$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new row (row 7) with the AI-related entry.
# Values are entered in this order so that new shared-string entries
# line up with the target workbook's shared string table ordering.
$ws.Range("A7").Value = "Company F"
$ws.Range("B7").Value = "Something Else"

$ws.Range("D7").Value = "ai_email@gmail.com"
$ws.Range("E7").Value = "file_for_ai.txt"
$ws.Range("F7").Value = "rayanakhtar1203@gmail.com"
$ws.Range("G7").Value = "AI internship"

$ws.Range("C7").Value = "ai_template_1.txt"

# Add hyperlinks to the recruiter-email columns (D and F), matching the
# pattern used by the existing rows, then re-apply the built-in
# "Hyperlink" cell style so the formatting matches the other rows.
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:ai_email@gmail.com")
$ws.Range("D7").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F7"), "mailto:rayanakhtar1203@gmail.com")
$ws.Range("F7").Style = "Hyperlink"

# Update the active selection to C1, matching the saved workbook view.
$ws.Range("C1").Select()
